$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Status text change ("Ready for handoff" -> "Handed back: in sync with en-US")
# This string is shared by the Overview sheet (columns B/C) and by the
# per-locale sheets (column C). Updating every cell that currently holds the
# old text keeps the shared string content correct everywhere it is used.
# ---------------------------------------------------------------------------
$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    foreach ($row in 1..$used.Rows.Count) {
        foreach ($col in 1..$used.Columns.Count) {
            $cell = $used.Cells.Item($row, $col)
            if ($cell.Value() -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet ("zh-cn") - add Latest Target File / Latest Handback File
# hyperlinks + values, and refresh the Latest Handback DateTime column.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/acb166ec52fa278f417bf0d53228fb9652fa6f28/e2e/765e91e4-d19f-455d-941b-b3676c6e74aa.md", "", "", "765e91e4-d19f-455d-941b-b3676c6e74aa.md")
$wsZh.Range("F2").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/61c5769578522483236fa3844cb44983ded11563/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/765e91e4-d19f-455d-941b-b3676c6e74aa.2fefb2047c8adc3042945d23eb9408831761cd98.zh-cn.xlf", "", "", "765e91e4-d19f-455d-941b-b3676c6e74aa.2fefb2047c8adc3042945d23eb9408831761cd98.zh-cn.xlf")
$wsZh.Range("G2").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/acb166ec52fa278f417bf0d53228fb9652fa6f28/e2e/a3ae426a-3eae-451e-a516-b0c3bfd235b4.md", "", "", "a3ae426a-3eae-451e-a516-b0c3bfd235b4.md")
$wsZh.Range("F3").Style = "HyperLink"

$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/61c5769578522483236fa3844cb44983ded11563/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/a3ae426a-3eae-451e-a516-b0c3bfd235b4.732b3cfd6e8e282cd2a2ef6874379d2a8c15aaff.zh-cn.xlf", "", "", "a3ae426a-3eae-451e-a516-b0c3bfd235b4.732b3cfd6e8e282cd2a2ef6874379d2a8c15aaff.zh-cn.xlf")
$wsZh.Range("G3").Style = "HyperLink"

$wsZh.Range("H2").Value = "2016-03-11 10:20:37"
$wsZh.Range("H3").Value = "2016-03-11 10:20:37"

# ---------------------------------------------------------------------------
# de-de sheet ("de-de") - same additions as above, using the de-de targets.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/acb166ec52fa278f417bf0d53228fb9652fa6f28/e2e/765e91e4-d19f-455d-941b-b3676c6e74aa.md", "", "", "765e91e4-d19f-455d-941b-b3676c6e74aa.md")
$wsDe.Range("F2").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f492d4e9370a359ce8cbc814a502be605464c04/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/765e91e4-d19f-455d-941b-b3676c6e74aa.2fefb2047c8adc3042945d23eb9408831761cd98.de-de.xlf", "", "", "765e91e4-d19f-455d-941b-b3676c6e74aa.2fefb2047c8adc3042945d23eb9408831761cd98.de-de.xlf")
$wsDe.Range("G2").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/acb166ec52fa278f417bf0d53228fb9652fa6f28/e2e/a3ae426a-3eae-451e-a516-b0c3bfd235b4.md", "", "", "a3ae426a-3eae-451e-a516-b0c3bfd235b4.md")
$wsDe.Range("F3").Style = "HyperLink"

$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9f492d4e9370a359ce8cbc814a502be605464c04/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/a3ae426a-3eae-451e-a516-b0c3bfd235b4.732b3cfd6e8e282cd2a2ef6874379d2a8c15aaff.de-de.xlf", "", "", "a3ae426a-3eae-451e-a516-b0c3bfd235b4.732b3cfd6e8e282cd2a2ef6874379d2a8c15aaff.de-de.xlf")
$wsDe.Range("G3").Style = "HyperLink"

$wsDe.Range("H2").Value = "2016-03-11 10:20:43"
$wsDe.Range("H3").Value = "2016-03-11 10:20:43"
